$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-01 Saturday", "2024-06-02 Sunday"),
    @("45÷8=5, 5", "30÷2=15, 0"),
    @("33÷9=3, 6", "85÷9=9, 4"),
    @("92÷4=23, 0", "42÷2=21, 0"),
    @("60÷2=30, 0", "98÷2=49, 0"),
    @("14÷8=1, 6", "95÷2=47, 1"),
    @("50÷2=25, 0", "10÷5=2, 0"),
    @("88÷9=9, 7", "41÷2=20, 1"),
    @("99÷4=24, 3", "50÷9=5, 5"),
    @("34÷8=4, 2", "65÷6=10, 5"),
    @("76÷4=19, 0", "76÷3=25, 1"),
    @("84÷9=9, 3", "30÷4=7, 2"),
    @("40÷7=5, 5", "35÷3=11, 2"),
    @("81÷2=40, 1", "88÷6=14, 4"),
    @("26÷9=2, 8", "71÷8=8, 7"),
    @("11÷5=2, 1", "56÷4=14, 0"),
    @("66÷7=9, 3", "31÷5=6, 1"),
    @("68÷2=34, 0", "71÷6=11, 5"),
    @("60÷8=7, 4", "54÷4=13, 2"),
    @("76÷9=8, 4", "53÷5=10, 3"),
    @("90÷2=45, 0", "87÷2=43, 1"),
    @("17÷5=3, 2", "63÷9=7, 0"),
    @("92÷5=18, 2", "85÷6=14, 1"),
    @("13÷9=1, 4", "42÷2=21, 0"),
    @("93÷6=15, 3", "53÷4=13, 1"),
    @("32÷3=10, 2", "18÷5=3, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
